$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 82, shifting existing rows 82:165 down to 83:166
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new record
$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = 44741
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 100112013
$ws.Range("G82").Value = "Alcachofa"
$ws.Range("H82").Value = "Española"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 50
$ws.Range("K82").Value = 22000
$ws.Range("L82").Value = 22000
$ws.Range("M82").Value = 22000
$ws.Range("N82").Value = "$/caja 30 unidades"
$ws.Range("O82").Value = "Provincia de Limarí"
$ws.Range("P82").Value = 733
$ws.Range("Q82").Value = 30
$ws.Range("R82").Value = "Hortaliza"
